$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.595171257880907
$ws.Range("C2").Value = 0.1997189457370325
$ws.Range("D2").Value = 0.1640302068917023
$ws.Range("F2").Value = 1.427060993182451
$ws.Range("G2").Value = 0.002465091533515773
$ws.Range("I2").Value = 0.835027060010944
$ws.Range("J2").Value = 0.1679191627495049
$ws.Range("L2").Value = 0.4220987939883969
$ws.Range("O2").Value = 3.483305048085299
$ws.Range("B3").Value = 1.458061052358687
$ws.Range("C3").Value = 0.1810641220763216
$ws.Range("D3").Value = 0.163052254729827
$ws.Range("F3").Value = 1.437753617176959
$ws.Range("G3").Value = 0.002468004936192551
$ws.Range("I3").Value = 0.8482709281717433
$ws.Range("J3").Value = 0.1699010648156127
$ws.Range("L3").Value = 0.4117511432170602
$ws.Range("O3").Value = 3.516362382676448
$ws.Range("B4").Value = 1.373851724902124
$ws.Range("C4").Value = 0.1695529483934308
$ws.Range("D4").Value = 0.1624951693002075
$ws.Range("F4").Value = 1.445229784162997
$ws.Range("G4").Value = 0.002469889584334751
$ws.Range("I4").Value = 0.8569684122360322
$ws.Range("J4").Value = 0.1711883932998424
$ws.Range("L4").Value = 0.4055297817785544
$ws.Range("O4").Value = 3.539096566826629
$ws.Range("B5").Value = 1.339532264368984
$ws.Range("C5").Value = 0.1648480123128877
$ws.Range("D5").Value = 0.16227911237943
$ws.Range("F5").Value = 1.44850530764149
$ws.Range("G5").Value = 0.002470681757301968
$ws.Range("I5").Value = 0.8606548052805323
$ws.Range("J5").Value = 0.1717307120527201
$ws.Range("L5").Value = 0.4030279691878604
$ws.Range("O5").Value = 3.548973061694937
$ws.Range("B6").Value = 1.333833403282199
$ws.Range("C6").Value = 0.1640659222935028
$ws.Range("D6").Value = 0.1622438997601918
$ws.Range("F6").Value = 1.449063028890954
$ws.Range("G6").Value = 0.002470814758287025
$ws.Range("I6").Value = 0.8612755059535253
$ws.Range("J6").Value = 0.1718218343568552
$ws.Range("L6").Value = 0.402614571340564
$ws.Range("O6").Value = 3.550649998675851
$ws.Range("B7").Value = 1.373388891619754
$ws.Range("C7").Value = 0.1694895524206856
$ws.Range("D7").Value = 0.1624922110355982
$ws.Range("F7").Value = 1.445273032267643
$ws.Range("G7").Value = 0.002469900170034108
$ws.Range("I7").Value = 0.8570175530759379
$ws.Range("J7").Value = 0.1711956354300241
$ws.Range("L7").Value = 0.4054959057633027
$ws.Range("O7").Value = 3.539227287048561
$ws.Range("B8").Value = 1.547901806644745
$ws.Range("C8").Value = 0.1932987758278841
$ws.Range("D8").Value = 0.1636840354335476
$ws.Range("F8").Value = 1.430558691328805
$ws.Range("G8").Value = 0.002466076233277835
$ws.Range("I8").Value = 0.83947604615609
$ws.Range("J8").Value = 0.1685879111479522
$ws.Range("L8").Value = 0.4185036237979602
$ws.Range("O8").Value = 3.49419706683716
$ws.Range("B9").Value = 1.889852739295748
$ws.Range("C9").Value = 0.239524987628414
$ws.Range("D9").Value = 0.1663633795384172
$ws.Range("F9").Value = 1.408937147651827
$ws.Range("G9").Value = 0.002459334408473368
$ws.Range("I9").Value = 0.8095721666738775
$ws.Range("J9").Value = 0.1640324636751069
$ws.Range("L9").Value = 0.4450523290965265
$ws.Range("O9").Value = 3.425258083846984
$ws.Range("B10").Value = 2.140832052283429
$ws.Range("C10").Value = 0.2731932537968476
$ws.Range("D10").Value = 0.1685381015435254
$ws.Range("F10").Value = 1.397470834325233
$ws.Range("G10").Value = 0.002454838011865277
$ws.Range("I10").Value = 0.7903503664421905
$ws.Range("J10").Value = 0.1610250686445207
$ws.Range("L10").Value = 0.4651838954137872
$ws.Range("O10").Value = 3.386457225768538
$ws.Range("B11").Value = 2.254937450459977
$ws.Range("C11").Value = 0.2884437591832523
$ws.Range("D11").Value = 0.1695717321124377
$ws.Range("F11").Value = 1.393216176383319
$ws.Range("G11").Value = 0.002452890698961648
$ws.Range("I11").Value = 0.7822044726774173
$ws.Range("J11").Value = 0.1597304599549201
$ws.Range("L11").Value = 0.4744766859455751
$ws.Range("O11").Value = 3.371387920683986
$ws.Range("B12").Value = 2.298134709530359
$ws.Range("C12").Value = 0.2942090779358466
$ws.Range("D12").Value = 0.1699694704335997
$ws.Range("F12").Value = 1.391743457384692
$ws.Range("G12").Value = 0.002452167337596163
$ws.Range("I12").Value = 0.7792060197177122
$ws.Range("J12").Value = 0.1592507790292181
$ws.Range("L12").Value = 0.4780148337391239
$ws.Range("O12").Value = 3.366053517487018
$ws.Range("B13").Value = 2.288831978229268
$ws.Range("C13").Value = 0.2929678505733762
$ws.Range("D13").Value = 0.1698835298786463
$ws.Range("F13").Value = 1.392054474378767
$ws.Range("G13").Value = 0.002452322502923463
$ws.Range("I13").Value = 0.7798479523647757
$ws.Range("J13").Value = 0.1593536173608587
$ws.Range("L13").Value = 0.4772519812200784
$ws.Range("O13").Value = 3.367185818250647
$ws.Range("B14").Value = 2.258491568918544
$ws.Range("C14").Value = 0.2889182719755468
$ws.Range("D14").Value = 0.1696043278159038
$ws.Range("F14").Value = 1.393092239606531
$ws.Range("G14").Value = 0.002452830906399179
$ws.Range("I14").Value = 0.7819560590225194
$ws.Range("J14").Value = 0.1596907847245923
$ws.Range("L14").Value = 0.4747673885606218
$ws.Range("O14").Value = 3.37094159572959
$ws.Range("B15").Value = 2.239905560191573
$ws.Range("C15").Value = 0.2864365133674482
$ws.Range("D15").Value = 0.1694341306640794
$ws.Range("F15").Value = 1.393745933333406
$ws.Range("G15").Value = 0.002453144145275267
$ws.Range("I15").Value = 0.7832585695648966
$ws.Range("J15").Value = 0.1598986844353716
$ws.Range("L15").Value = 0.4732479937657672
$ws.Range("O15").Value = 3.373290589373397
$ws.Range("B16").Value = 2.133373452131991
$ws.Range("C16").Value = 0.2721952549577509
$ws.Range("D16").Value = 0.1684714391970061
$ws.Range("F16").Value = 1.397768247623134
$ws.Range("G16").Value = 0.002454967242150585
$ws.Range("I16").Value = 0.7908947741517132
$ws.Range("J16").Value = 0.1611111522921735
$ws.Range("L16").Value = 0.4645792841762244
$ws.Range("O16").Value = 3.387494054275578
$ws.Range("B17").Value = 2.068000718956569
$ws.Range("C17").Value = 0.2634417370171889
$ws.Range("D17").Value = 0.1678921789841681
$ws.Range("F17").Value = 1.400482176774581
$ws.Range("G17").Value = 0.002456110736954625
$ws.Range("I17").Value = 0.7957326938775999
$ws.Range("J17").Value = 0.1618737755121087
$ws.Range("L17").Value = 0.4592956903405678
$ws.Range("O17").Value = 3.396869158290031
$ws.Range("B18").Value = 2.030393950979146
$ws.Range("C18").Value = 0.2584008120596764
$ws.Range("D18").Value = 0.1675631802227215
$ws.Range("F18").Value = 1.4021336405212
$ws.Range("G18").Value = 0.002456777684704327
$ws.Range("I18").Value = 0.7985716307850161
$ws.Range("J18").Value = 0.162319331930191
$ws.Range("L18").Value = 0.4562694117735191
$ws.Range("O18").Value = 3.402504395305357
$ws.Range("B19").Value = 2.017659972326612
$ws.Range("C19").Value = 0.2566929978945325
$ws.Range("D19").Value = 0.1674525057530829
$ws.Range("F19").Value = 1.402708333021586
$ws.Range("G19").Value = 0.002457005090459951
$ws.Range("I19").Value = 0.7995425084036079
$ws.Range("J19").Value = 0.1624713776322686
$ws.Range("L19").Value = 0.4552469546497093
$ws.Range("O19").Value = 3.404454087529331
$ws.Range("B20").Value = 2.074960407152332
$ws.Range("C20").Value = 0.2643742009987022
$ws.Range("D20").Value = 0.1679534103527658
$ws.Range("F20").Value = 1.400183908403186
$ws.Range("G20").Value = 0.00245598805444242
$ws.Range("I20").Value = 0.7952118614535095
$ws.Range("J20").Value = 0.1617918772596205
$ws.Range("L20").Value = 0.4598568246183135
$ws.Range("O20").Value = 3.395846014710912
$ws.Range("B21").Value = 2.267403627015199
$ws.Range("C21").Value = 0.290107997601865
$ws.Range("D21").Value = 0.169686164979673
$ws.Range("F21").Value = 1.392783664522184
$ws.Range("G21").Value = 0.002452681195367833
$ws.Range("I21").Value = 0.7813345154124107
$ws.Range("J21").Value = 0.159591463965671
$ws.Range("L21").Value = 0.4754966551400912
$ws.Range("O21").Value = 3.36982832998072
$ws.Range("B22").Value = 2.393105458725245
$ws.Range("C22").Value = 0.3068697397294216
$ws.Range("D22").Value = 0.1708554629339716
$ws.Range("F22").Value = 1.388754120236612
$ws.Range("G22").Value = 0.002450601800184853
$ws.Range("I22").Value = 0.7727675175353994
$ws.Range("J22").Value = 0.1582149117439577
$ws.Range("L22").Value = 0.4858298229463998
$ws.Range("O22").Value = 3.354993035576001
$ws.Range("B23").Value = 2.326023315896464
$ws.Range("C23").Value = 0.2979289823560975
$ws.Range("D23").Value = 0.1702280322792546
$ws.Range("F23").Value = 1.390830870754314
$ws.Range("G23").Value = 0.002451704146835154
$ws.Range("I23").Value = 0.7772938270656731
$ws.Range("J23").Value = 0.1589439740343588
$ws.Range("L23").Value = 0.4803046746495596
$ws.Range("O23").Value = 3.362712201260308
$ws.Range("B24").Value = 2.071814003929035
$ws.Range("C24").Value = 0.2639526602094122
$ws.Range("D24").Value = 0.1679257151076925
$ws.Range("F24").Value = 1.400318471494799
$ws.Range("G24").Value = 0.002456043489524608
$ws.Range("I24").Value = 0.7954471505328655
$ws.Range("J24").Value = 0.1618288812990469
$ws.Range("L24").Value = 0.4596031005051202
$ws.Range("O24").Value = 3.396307813679584
$ws.Range("B25").Value = 1.797383912054897
$ws.Range("C25").Value = 0.2270704065113591
$ws.Range("D25").Value = 0.1656021618207362
$ws.Range("F25").Value = 1.414010998755032
$ws.Range("G25").Value = 0.00246107769695477
$ws.Range("I25").Value = 0.8171798800682808
$ws.Range("J25").Value = 0.1652051530357141
$ws.Range("L25").Value = 0.4377596443673184
$ws.Range("O25").Value = 3.441830266865168
